$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 329 (pushes old rows 329-334 down to 331-336,
# carrying their formatting - column D's date style - along with them).
$ws.Rows("329:330").Insert()

# New row 329: "Pintón" quality entry for the week of 2021-09-09 (serial 44448).
$ws.Range("A329").Value = 7
$ws.Range("B329").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C329").Value = "Ñuble"
$ws.Range("D329").Value = 44448
$ws.Range("E329").Value = 16
$ws.Range("F329").Value = "Fruta"
$ws.Range("G329").Value = 100108
$ws.Range("H329").Value = "Tropicales y subtropicales"
$ws.Range("I329").Value = 100108006
$ws.Range("J329").Value = "Plátano"
$ws.Range("K329").Value = "Sin especificar"
$ws.Range("L329").Value = "Pintón"
$ws.Range("M329").Value = 180
$ws.Range("N329").Value = 19000
$ws.Range("O329").Value = 19000
$ws.Range("P329").Value = 19000
$ws.Range("Q329").Value = "$/caja 20 kilos"
$ws.Range("R329").Value = "Ecuador"
$ws.Range("S329").Value = 950
$ws.Range("T329").Value = 20

# New row 330: "Primera Pintón" quality entry for the same week.
$ws.Range("A330").Value = 7
$ws.Range("B330").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C330").Value = "Ñuble"
$ws.Range("D330").Value = 44448
$ws.Range("E330").Value = 16
$ws.Range("F330").Value = "Fruta"
$ws.Range("G330").Value = 100108
$ws.Range("H330").Value = "Tropicales y subtropicales"
$ws.Range("I330").Value = 100108006
$ws.Range("J330").Value = "Plátano"
$ws.Range("K330").Value = "Sin especificar"
$ws.Range("L330").Value = "Primera Pintón"
$ws.Range("M330").Value = 400
$ws.Range("N330").Value = 20000
$ws.Range("O330").Value = 21000
$ws.Range("P330").Value = 20500
$ws.Range("Q330").Value = "$/caja 20 kilos"
$ws.Range("R330").Value = "Ecuador"
$ws.Range("S330").Value = 1025
$ws.Range("T330").Value = 20
